$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Text values -------------------------------------------------------
# New shared strings must be created in this exact order so they land at
# sharedStrings indices 30-36 in the same sequence as the target workbook:
#   30 Comentarios
#   31 consultas
#   32 el formulario de consulta funciona correctamente
#   33 visualizacion de consultas
#   34 las consultas solo se ven en el perfil de un superuser donde
#      aparecen los datos que envio la persona en el formulario
#   35 Allen/Rio Negro
#   36 los comentarios de cargan correctamente y se visualizan solo en el
#      curso donde se realizo
$ws.Range("B29").Value = "Comentarios"
$ws.Range("B30").Value = "consultas"
$ws.Range("C30").Value = "el formulario de consulta funciona correctamente"
$ws.Range("B31").Value = "visualizacion de consultas"
$ws.Range("C31").Value = "las consultas solo se ven en el perfil de un superuser donde aparecen los datos que envio la persona en el formulario"
$ws.Range("D26").Value = "Allen/Rio Negro"
$ws.Range("C29").Value = "los comentarios de cargan correctamente y se visualizan solo en el curso donde se realizo"

# --- Case numbers --------------------------------------------------------
$ws.Range("A29").Value = 1
$ws.Range("A30").Value = 2
$ws.Range("A31").Value = 3

# --- Approved column -------------------------------------------------------
$ws.Range("G29").Value = "ok"
$ws.Range("G30").Value = "ok"
$ws.Range("G31").Value = "ok"

# --- Dates -----------------------------------------------------------------
# Set the raw date serials first, then apply a short-date display format.
# Formatting D29 first creates one new cell style (numFmtId 14, thin border);
# copying that formatting onto D30:D31 re-uses the very same style index
# instead of minting a near-duplicate style per cell.
$ws.Range("D29").Value = 45217
$ws.Range("D30").Value = 45218
$ws.Range("D31").Value = 45219

$ws.Range("D29").NumberFormat = "mm-dd-yy"
$ws.Range("D29").Copy()
$ws.Range("D30:D31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection ---------------------------------------------------------
$ws.Range("C29").Select()
